$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 32   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/8/2025  Through  9/14/2025"

# --- Crime data table updates (rows 15-31) ---
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("I15").Value = 7
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -12.5
$ws.Range("M15").Value = -12.5
$ws.Range("N15").Value = -61.111111111111
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 69
$ws.Range("J16").Value = 71
$ws.Range("K16").Value = -2.81690140845
$ws.Range("L16").Value = 11.290322580645
$ws.Range("M16").Value = 4.545454545454
$ws.Range("N16").Value = -83.764705882352
$ws.Range("C17").Value = 3
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -46.153846153846
$ws.Range("I17").Value = 55
$ws.Range("J17").Value = 87
$ws.Range("K17").Value = -36.781609195402
$ws.Range("L17").Value = -31.25
$ws.Range("M17").Value = 19.565217391304
$ws.Range("N17").Value = -19.117647058823
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 96
$ws.Range("J18").Value = 95
$ws.Range("K18").Value = 1.052631578947
$ws.Range("L18").Value = 15.662650602409
$ws.Range("M18").Value = 26.315789473684
$ws.Range("N18").Value = -84.834123222748
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = 1.851851851851
$ws.Range("I19").Value = 517
$ws.Range("J19").Value = 543
$ws.Range("K19").Value = -4.788213627992
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 12.636165577342
$ws.Range("N19").Value = -62.013225569434
$ws.Range("D20").Value = 5
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -70
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = -17.5
$ws.Range("L20").Value = -55.405405405405
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -95.864661654135
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -40.74074074074
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -5.617977528089
$ws.Range("I21").Value = 777
$ws.Range("J21").Value = 843
$ws.Range("K21").Value = -7.829181494661
$ws.Range("L21").Value = -6.159420289855
$ws.Range("M21").Value = 14.601769911504
$ws.Range("N21").Value = -76.483050847457
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -50
$ws.Range("D23").Value = 3
$ws.Range("G23").Value = 7
$ws.Range("J23").Value = 32
$ws.Range("K23").Value = -50
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -21.875
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = 3.960396039603
$ws.Range("I24").Value = 1019
$ws.Range("J24").Value = 907
$ws.Range("K24").Value = 12.348401323043
$ws.Range("L24").Value = 23.665048543689
$ws.Range("M24").Value = 35.505319148936
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -26.086956521739
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 72
$ws.Range("H25").Value = -16.666666666666
$ws.Range("I25").Value = 751
$ws.Range("J25").Value = 710
$ws.Range("K25").Value = 5.774647887323
$ws.Range("L25").Value = 19.96805111821
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 168
$ws.Range("J26").Value = 165
$ws.Range("K26").Value = 1.818181818181
$ws.Range("L26").Value = 4.347826086956
$ws.Range("M26").Value = -13.40206185567
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("I27").Value = 10
$ws.Range("K27").Value = 11.111111111111
$ws.Range("L27").Value = -37.5
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 0
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 1
$ws.Range("G28").NumberFormat = "#,##0"
$ws.Range("H28").Value = 300
$ws.Range("H28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I28").Value = 31
$ws.Range("J28").Value = 19
$ws.Range("K28").Value = 63.157894736842
$ws.Range("L28").Value = 10.714285714285
$ws.Range("D29").Value = 2
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G29").Value = 2
$ws.Range("G29").NumberFormat = "#,##0"
$ws.Range("H29").Value = -100
$ws.Range("H29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J29").Value = 4
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J30").Value = 3
$ws.Range("F31").Value = 1
$ws.Range("F31").NumberFormat = "#,##0"
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 6
$ws.Range("K31").Value = -40
$ws.Range("L31").Value = -25
